# "shear wall work continue"
# - Width/Depth footprint values drop from 120 ft to 90 ft
# - The "Width" label (B4) loses its stray bold formatting (matches B5/B6)
# - Active selection moves to C4 (Width's value cell)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C4").Value = 90
$ws.Range("C5").Value = 90

# B4 ("Width") had an orphaned "applyFont" style left over from an earlier
# bold toggle; unbolding it brings it back in line with the plain B5/B6 cells.
$ws.Range("B4").Font.Bold = $false

$ws.Range("C4").Select()
